$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (removed rowId/isHeader columns concept -> generic sheet name)
$ws.Name = "RelationShip"

# New row 2: first component code being tracked
$ws.Range("A2").Value = "org.emoflon.ibex.tgg.core.language"

# Highlight the new row with the "Orange, Accent 6, Lighter 40%" theme fill.
# Apply to a single cell first so the intermediate style objects created while
# building the themed fill don't linger as used cell formats, then replicate
# that exact format across the rest of the row.
$a2 = $ws.Range("A2")
$a2.Interior.ThemeColor = 10
$a2.Interior.TintAndShade = 0.39997558519241921

$a2.Copy()
$ws.Range("B2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Select column C, matching the author's click-through while reviewing data
[void]$ws.Columns("C").Select()
